$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Data for each sheet: header labels (row1, columns B/C) + row label/value
# pairs placed starting at row 2 (col A = label, col B/C = distilgpt2/gpt2).
# ---------------------------------------------------------------------------
$sheetsData = @(
    @{
        Name = "Architecture"
        Rows = @(
            @("parameters", 81912576, 124439808),
            @("layers", 86, 164),
            @("vocab_size", 50257, 50257),
            @("hidden_size", 768, 768),
            @("num_attention_heads", 12, 12)
        )
    },
    @{
        Name = "Speed"
        Rows = @(
            @("avg_inference_time", 5.574, 9.651),
            @("tokens_per_second", 8.97, 5.181),
            @("std_inference_time", 0.037, 0.038)
        )
    },
    @{
        Name = "Quality"
        Rows = @(
            @("lexical_diversity", 0.211, 0.235),
            @("avg_length", 43.667, 64.667),
            @("repetition_rate", 1, 0.654)
        )
    },
    @{
        Name = "Memory"
        Rows = @(
            @("model_size_mb", 156.24, 237.35),
            @("vocab_size_mb", 0.1, 0.1)
        )
    }
)

# Build one shared "bold / thin-box-border / centered-top" template style and
# remember which ranges (per sheet) need it -- this way the style is created
# only once (single new font + single new border + single new cellXf) and
# just copy/paste-special'd onto every sheet, instead of re-deriving it per
# sheet (which would otherwise mint a fresh cellXf each time).
$firstSheet = $wb.Worksheets.Item($sheetsData[0].Name)
$template = $firstSheet.Range("Z100")
$template.Font.Bold = $true
$template.Borders.LineStyle = 1
$template.HorizontalAlignment = -4108
$template.VerticalAlignment = -4160
$template.Copy()

foreach ($sheetData in $sheetsData) {
    $ws = $wb.Worksheets.Item($sheetData.Name)

    # --- Header row -----------------------------------------------------
    $ws.Range("B1").Value = "distilgpt2"
    $ws.Range("C1").Value = "gpt2"

    # --- Body rows --------------------------------------------------------
    $r = 2
    foreach ($row in $sheetData.Rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $r = $r + 1
    }
    $lastRow = $r - 1

    # --- Formatting: bold, centered/top-aligned, thin box border, applied
    # only to the label cells (header row B/C, and column A body labels) ---
    $ws.Range("B1:C1").PasteSpecial(-4122)
    $ws.Range("A2:A" + $lastRow).PasteSpecial(-4122)
}

$template.Clear()

Write-Host "done"
